$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 10488.833
$ws.Range("J17").Value = 10488.833
$ws.Range("L17").Value = 31466.499
$ws.Range("N17").Value = -31802.499
$ws.Range("H88").Value = 17370.889
$ws.Range("I88").Value = 7697.923
$ws.Range("K88").Value = 7697.923
$ws.Range("M88").Value = -7291.923
$ws.Range("H91").Value = 17370.889
$ws.Range("I91").Value = 7697.923
$ws.Range("K91").Value = 7697.923
$ws.Range("M91").Value = -6293.923
$ws.Range("H112").Value = 1263.125
$ws.Range("J112").Value = 1271.6129
$ws.Range("L112").Value = 3814.8387
$ws.Range("N112").Value = -6030.8387
$ws.Range("H129").Value = 355644.03
$ws.Range("I129").Value = 2403199.5
$ws.Range("J129").Value = 2617.2415
$ws.Range("K129").Value = 7209598.5
$ws.Range("L129").Value = 7851.7245
$ws.Range("M129").Value = -7204598.5
$ws.Range("N129").Value = -17851.7245
$ws.Range("H138").Value = 2880.602
$ws.Range("I138").Value = 2254.5518
$ws.Range("J138").Value = 3164.2812
$ws.Range("K138").Value = 6763.655400000001
$ws.Range("L138").Value = 9492.8436
$ws.Range("M138").Value = -1623.655400000001
$ws.Range("N138").Value = -19772.8436

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22968.236
$ws.Range("I32").Value = 22774.693
$ws.Range("J32").Value = 40000
$ws.Range("K32").Value = 22774.693
$ws.Range("L32").Value = 40000
$ws.Range("M32").Value = -22487.693
$ws.Range("N32").Value = -40574
$ws.Range("H61").Value = 1630.9701
$ws.Range("I61").Value = 1431.5
$ws.Range("J61").Value = 1966.08
$ws.Range("K61").Value = 1431.5
$ws.Range("L61").Value = 1966.08
$ws.Range("M61").Value = -1219.5
$ws.Range("N61").Value = -2390.08
$ws.Range("H74").Value = 2215.4062
$ws.Range("I74").Value = 1794.7142
$ws.Range("J74").Value = 3018.5454
$ws.Range("K74").Value = 1794.7142
$ws.Range("L74").Value = 3018.5454
$ws.Range("M74").Value = -920.7141999999999
$ws.Range("N74").Value = -4766.5454
$ws.Range("H77").Value = 2215.4062
$ws.Range("I77").Value = 1794.7142
$ws.Range("J77").Value = 3018.5454
$ws.Range("K77").Value = 8973.571
$ws.Range("L77").Value = 15092.727
$ws.Range("M77").Value = -4605.571
$ws.Range("N77").Value = -23828.727
$ws.Range("H110").Value = 2090.8438
$ws.Range("I110").Value = 2207.3809
$ws.Range("J110").Value = 1868.3636
$ws.Range("K110").Value = 2207.3809
$ws.Range("L110").Value = 1868.3636
$ws.Range("M110").Value = -162.3809000000001
$ws.Range("N110").Value = -5958.3636
$ws.Range("H132").Value = 7464818
$ws.Range("I132").Value = 14707316
$ws.Range("J132").Value = 2849.818
$ws.Range("K132").Value = 44121948
$ws.Range("L132").Value = 8549.454000000002
$ws.Range("M132").Value = -44119418
$ws.Range("N132").Value = -13609.454
$ws.Range("H136").Value = 1630.9701
$ws.Range("I136").Value = 1431.5
$ws.Range("J136").Value = 1966.08
$ws.Range("K136").Value = 4294.5
$ws.Range("L136").Value = 5898.24
$ws.Range("M136").Value = -1744.5
$ws.Range("N136").Value = -10998.24

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2216.4138
$ws.Range("I105").Value = 2120.3845
$ws.Range("J105").Value = 2294.4375
$ws.Range("K105").Value = 2120.3845
$ws.Range("L105").Value = 2294.4375
$ws.Range("M105").Value = -373.3845000000001
$ws.Range("N105").Value = -5788.4375

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16420.375
$ws.Range("I31").Value = 4180.75
$ws.Range("J31").Value = 28660
$ws.Range("K31").Value = 4180.75
$ws.Range("L31").Value = 28660
$ws.Range("M31").Value = -3885.75
$ws.Range("N31").Value = -29250
$ws.Range("H34").Value = 16420.375
$ws.Range("I34").Value = 4180.75
$ws.Range("J34").Value = 28660
$ws.Range("K34").Value = 4180.75
$ws.Range("L34").Value = 28660
$ws.Range("M34").Value = -3978.75
$ws.Range("N34").Value = -29064
$ws.Range("H132").Value = 24781.344
$ws.Range("I132").Value = 1296.341
$ws.Range("J132").Value = 85566.06
$ws.Range("K132").Value = 3889.023
$ws.Range("L132").Value = 256698.18
$ws.Range("M132").Value = -1359.023
$ws.Range("N132").Value = -261758.18

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 570.3
$ws.Range("I47").Value = 240.6
$ws.Range("K47").Value = 721.8
$ws.Range("M47").Value = -290.8

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 37083
$ws.Range("J92").Value = 37083
$ws.Range("L92").Value = 37083
$ws.Range("N92").Value = -40827
$ws.Range("H122").Value = 1640
$ws.Range("I122").Value = 1575
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 4725
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -2275
$ws.Range("N122").Value = -10600
$ws.Range("H132").Value = 3316
$ws.Range("I132").Value = 2471.4375
$ws.Range("J132").Value = 4442.0835
$ws.Range("K132").Value = 7414.3125
$ws.Range("L132").Value = 13326.2505
$ws.Range("M132").Value = -4884.3125
$ws.Range("N132").Value = -18386.2505

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1950.5
$ws.Range("I100").Value = 1764.1818
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 1764.1818
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -1223.1818
$ws.Range("N100").Value = -5082
$ws.Range("H119").Value = 46664
$ws.Range("J119").Value = 46664
$ws.Range("L119").Value = 46664
$ws.Range("N119").Value = -56340
$ws.Range("H123").Value = 40311
$ws.Range("J123").Value = 40311
$ws.Range("L123").Value = 40311
$ws.Range("N123").Value = -50111

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 74497
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 74497
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 74497
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -74959
$ws.Range("H100").Value = 83333864
$ws.Range("I100").Value = 640.8
$ws.Range("J100").Value = 500000000
$ws.Range("K100").Value = 1281.6
$ws.Range("L100").Value = 1000000000
$ws.Range("M100").Value = -740.5999999999999
$ws.Range("N100").Value = -1000001082
$ws.Range("H132").Value = 1499.2041
$ws.Range("I132").Value = 1147.7317
$ws.Range("J132").Value = 3300.5
$ws.Range("K132").Value = 3443.1951
$ws.Range("L132").Value = 9901.5
$ws.Range("M132").Value = -913.1950999999999
$ws.Range("N132").Value = -14961.5
$ws.Range("H134").Value = 74497
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 74497
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 223491
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -228561
